# edit.ps1 - Update match-odds values in the FlashScore weekly odds workbook
# (commit: "Atualizando o arquivo XLSX" - refreshed odds snapshot for 2024-10-17).
# Applies the per-cell numeric updates captured in the target diff, row by row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("M2").Value = 1.06
$ws.Range("O2").Value = 1.33
$ws.Range("V2").Value = 1.75

# Row 3
$ws.Range("M3").Value = 1.1
$ws.Range("O3").Value = 1.44
$ws.Range("P3").Value = 2.75
$ws.Range("V3").Value = 1.62

# Row 4
$ws.Range("G4").Value = 1.7
$ws.Range("H4").Value = 3.3
$ws.Range("I4").Value = 6
$ws.Range("K4").Value = 2.05
$ws.Range("L4").Value = 6
$ws.Range("M4").Value = 1.08
$ws.Range("N4").Value = 8
$ws.Range("O4").Value = 1.44
$ws.Range("P4").Value = 2.63
$ws.Range("Q4").Value = 2.35
$ws.Range("R4").Value = 1.57
$ws.Range("S4").Value = 1.5
$ws.Range("T4").Value = 2.5
$ws.Range("U4").Value = 2.2
$ws.Range("V4").Value = 1.62
$ws.Range("W4").Value = 5.5
$ws.Range("Z4").Value = 12
$ws.Range("AA4").Value = 17
$ws.Range("AC4").Value = 7
$ws.Range("AE4").Value = 21
$ws.Range("AF4").Value = 81
$ws.Range("AK4").Value = 67
$ws.Range("AL4").Value = 51
$ws.Range("AQ4").Value = 34
$ws.Range("AT4").Value = 2.5
$ws.Range("AU4").Value = 9.5
$ws.Range("AV4").Value = 81
$ws.Range("AW4").Value = 7
$ws.Range("AX4").Value = 34
$ws.Range("AZ4").Value = 126

# Row 5
$ws.Range("G5").Value = 2.88
$ws.Range("I5").Value = 2.45
$ws.Range("M5").Value = 1.06
$ws.Range("N5").Value = 9.5
$ws.Range("O5").Value = 1.33
$ws.Range("P5").Value = 3.25
$ws.Range("Q5").Value = 2.05
$ws.Range("R5").Value = 1.75
$ws.Range("S5").Value = 1.44
$ws.Range("T5").Value = 2.63
$ws.Range("U5").Value = 1.8
$ws.Range("V5").Value = 1.95
$ws.Range("W5").Value = 8.5
$ws.Range("X5").Value = 13
$ws.Range("AB5").Value = 34
$ws.Range("AC5").Value = 9.5
$ws.Range("AL5").Value = 21
$ws.Range("AN5").Value = 4.75
$ws.Range("AP5").Value = 26
$ws.Range("AR5").Value = 81
$ws.Range("AS5").Value = 201
$ws.Range("AT5").Value = 2.63

# Row 6
$ws.Range("G6").Value = 1.39
$ws.Range("H6").Value = 4.4
$ws.Range("I6").Value = 6.1
$ws.Range("J6").Value = 1.85
$ws.Range("K6").Value = 2.42
$ws.Range("L6").Value = 5.8
$ws.Range("O6").Value = 1.15
$ws.Range("P6").Value = 4
$ws.Range("U6").Value = 1.91
$ws.Range("V6").Value = 1.85
$ws.Range("W6").Value = 6.5
$ws.Range("X6").Value = 6
$ws.Range("Y6").Value = 7.1
$ws.Range("Z6").Value = 7.8
$ws.Range("AA6").Value = 9.25
$ws.Range("AB6").Value = 20
$ws.Range("AC6").Value = 13
$ws.Range("AD6").Value = 7.7
$ws.Range("AE6").Value = 15
$ws.Range("AF6").Value = 60
$ws.Range("AG6").Value = 400
$ws.Range("AH6").Value = 15
$ws.Range("AI6").Value = 32
$ws.Range("AJ6").Value = 16.5
$ws.Range("AK6").Value = 90
$ws.Range("AL6").Value = 50
$ws.Range("AM6").Value = 45
$ws.Range("AN6").Value = 3.3
$ws.Range("AO6").Value = 6.2
$ws.Range("AQ6").Value = 17
$ws.Range("AU6").Value = 8
$ws.Range("AV6").Value = 70
$ws.Range("AW6").Value = 7.8
$ws.Range("AX6").Value = 35
$ws.Range("AY6").Value = 35
$ws.Range("BA6").Value = 250
$ws.Range("BB6").Value = 450

# Row 7
$ws.Range("G7").Value = 2.82
$ws.Range("I7").Value = 2.4
$ws.Range("J7").Value = 3.3
$ws.Range("K7").Value = 2.07
$ws.Range("L7").Value = 3
$ws.Range("O7").Value = 1.23
$ws.Range("P7").Value = 3.35
$ws.Range("Q7").Value = 1.7
$ws.Range("R7").Value = 1.91
$ws.Range("U7").Value = 1.53
$ws.Range("V7").Value = 2.18
$ws.Range("W7").Value = 10.75
$ws.Range("X7").Value = 16.5
$ws.Range("Y7").Value = 10
$ws.Range("AB7").Value = 25
$ws.Range("AC7").Value = 11
$ws.Range("AD7").Value = 6.2
$ws.Range("AH7").Value = 9.5
$ws.Range("AI7").Value = 13.5
$ws.Range("AJ7").Value = 9
$ws.Range("AK7").Value = 27
$ws.Range("AL7").Value = 18.5
$ws.Range("AM7").Value = 23
$ws.Range("AN7").Value = 4.9
$ws.Range("AO7").Value = 15
$ws.Range("AP7").Value = 19.5
$ws.Range("AR7").Value = 90
$ws.Range("AU7").Value = 6.3
$ws.Range("AV7").Value = 50
$ws.Range("AW7").Value = 4.45
$ws.Range("AX7").Value = 13
$ws.Range("AY7").Value = 18.5
$ws.Range("AZ7").Value = 50
$ws.Range("BA7").Value = 75

# Row 10
$ws.Range("Q10").Value = 2.1
$ws.Range("R10").Value = 1.7
$ws.Range("U10").Value = 1.91
$ws.Range("V10").Value = 1.8

# Row 11
$ws.Range("V11").Value = 1.67

# Row 12
$ws.Range("V12").Value = 1.67

# Row 13
$ws.Range("V13").Value = 1.62
